$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "`"二级关键字是能够体现这段话的主旨、重点、观点或主要信息的词汇或短语。这里是针对每一个问询函问题提取关键字，且同一个问题可以提取多个不同的关键字标签。具体判断如下：
如果问题涉及下述关键字：三表科目名称
货币资金
交易性金融资产
衍生金融资产
应收票据
应收账款
应收款项融资
预付款项
应收股利
应收利息
其他应收款
存货
合同资产
持有待售资产
一年内到期的非流动资产
其他流动资产
债权投资
其他债权投资
其他权益工具投资
其他非流动金融资产
可供出售金融资产
持有至到期投资
长期应收款
长期股权投资
投资性房地产
固定资产
在建工程
生产性生物资产
油气资产
使用权资产
无形资产
开发支出
商誉
长期待摊费用
递延所得税资产
其他非流动资产
资产总计
短期借款
交易性金融负债
衍生金融负债
应付票据
应付账款
预收账款
合同负债
应付职工薪酬
应交税费
应付利息
应付股利
其他应付款
持有待售负债
一年内到期的非流动负债
其他流动负债
长期借款
应付债券
租赁负债
长期应付款
预计负债
递延所得税负债
递延收益
其他非流动负债
负债合计
实收资本(或股本)
其他权益工具
优先股
永续债
资本公积
库存股
其他综合收益
专项储备
盈余公积
未分配利润
外币报表折算差额
归属于母公司所有者权益（或股东权益）合计
少数股东权益
所有者权益（或股东权益）合计
营业总收入
营业收入
营业总成本
营业成本
税金及附加
销售费用
管理费用
研发费用
财务费用
利息费用
利息收入
其他收益
投资收益
汇兑收益
净敞口套期收益
公允价值变动收益
资产减值损失
信用减值损失
资产处置收益
营业利润
营业外收入
营业外支出
利润总额
所得税费用
净利润
持续经营净利润
终止经营净利润
少数股东损益
归属于母公司股东的净利润
其他综合收益
综合收益总额
归属于母公司所有者的综合收益总额
归属于少数股东的综合收益总额
每股收益
基本每股收益
稀释每股收益
销售商品、提供劳务收到的现金
收到的税费返还
收到其他与经营活动有关的现金
经营活动现金流入小计
购买商品、接受劳务支付的现金
支付给职工以及为职工支付的现金
支付的各项税费
支付其他与经营活动有关的现金
经营活动现金流出小计
经营活动产生的现金流量净额
收回投资收到的现金
取得投资收益收到的现金
处置固定资产、无形资产和其他长期资产收回的现金净额
处置子公司及其他营业单位收到的现金净额
收到其他与投资活动有关的现金
投资活动现金流入小计
购建固定资产、无形资产和其他长期资产支付的现金
投资支付的现金
取得子公司及其他营业单位支付的现金净额
支付的其他与投资活动有关的现金
投资活动现金流出小计
投资活动产生的现金流量净额
吸收投资收到的现金
取得借款收到的现金
收到其他与筹资活动有关的现金
发行债券收到的现金
筹资活动现金流入小计
偿还债务支付的现金
分配股利、利润或偿付利息支付的现金
支付其他与筹资活动有关的现金
筹资活动现金流出小计
筹资活动产生的现金流量净额
汇率变动对现金的影响
现金及现金等价物净增加额
期初现金及现金等价物余额
期末现金及现金等价物余额、偿债能力、毛利率、持续经营能力、研发水平、融资能力、内部控制、主要客户及供应商、业绩承诺、债务重组、股利分配、股权质押、资金占用、事务所变更、会计政策变更、会计估计变更、会计差错更正、关联交易、合同纠纷、违规担保。则打上相应的关键字标签。如果问询问题不涉及上述关键字，则打上”其他“这个标签。`""

$ws.Range("A2").Value = $newText
$ws.Range("A2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.5

$ws.Range("A2").Select()
